$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 63.666668
$ws.Range("I6").Value = 45.5
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 136.5
$ws.Range("L6").Value = 300
$ws.Range("M6").Value = -24.5
$ws.Range("N6").Value = -524

$ws.Range("H33").Value = 406.2381
$ws.Range("I33").Value = 364.78946
$ws.Range("K33").Value = 364.78946
$ws.Range("M33").Value = -135.78946

$ws.Range("H38").Value = 605.0769
$ws.Range("I38").Value = 227.66667
$ws.Range("J38").Value = 928.5714
$ws.Range("K38").Value = 683.00001
$ws.Range("L38").Value = 2785.7142
$ws.Range("M38").Value = -311.00001
$ws.Range("N38").Value = -3529.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14390.706
$ws.Range("I32").Value = 15041.473
$ws.Range("J32").Value = 10012.818
$ws.Range("K32").Value = 15041.473
$ws.Range("L32").Value = 10012.818
$ws.Range("M32").Value = -14754.473
$ws.Range("N32").Value = -10586.818

$ws.Range("H74").Value = 2134.6667
$ws.Range("I74").Value = 2174.75
$ws.Range("J74").Value = 2054.5
$ws.Range("K74").Value = 2174.75
$ws.Range("L74").Value = 2054.5
$ws.Range("M74").Value = -1300.75
$ws.Range("N74").Value = -3802.5

$ws.Range("H77").Value = 2134.6667
$ws.Range("I77").Value = 2174.75
$ws.Range("J77").Value = 2054.5
$ws.Range("K77").Value = 10873.75
$ws.Range("L77").Value = 10272.5
$ws.Range("M77").Value = -6505.75
$ws.Range("N77").Value = -19008.5

$ws.Range("H110").Value = 4320
$ws.Range("I110").Value = 4093.3333
$ws.Range("J110").Value = 5000
$ws.Range("K110").Value = 4093.3333
$ws.Range("L110").Value = 5000
$ws.Range("M110").Value = -2048.3333
$ws.Range("N110").Value = -9090

$ws.Range("H119").Value = 20633.334
$ws.Range("J119").Value = 20633.334
$ws.Range("L119").Value = 20633.334
$ws.Range("N119").Value = -30309.334

$ws.Range("H123").Value = 50429
$ws.Range("J123").Value = 50429
$ws.Range("L123").Value = 50429
$ws.Range("N123").Value = -60229

$ws.Range("H132").Value = 30062.684
$ws.Range("I132").Value = 3968.7693
$ws.Range("J132").Value = 86599.5
$ws.Range("K132").Value = 11906.3079
$ws.Range("L132").Value = 259798.5
$ws.Range("M132").Value = -9376.3079
$ws.Range("N132").Value = -264858.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26774.39
$ws.Range("I134").Value = 31188.742
$ws.Range("J134").Value = 1024
$ws.Range("K134").Value = 93566.226
$ws.Range("L134").Value = 3072
$ws.Range("M134").Value = -91031.226
$ws.Range("N134").Value = -8142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7466.7
$ws.Range("I31").Value = 12432.486
$ws.Range("J31").Value = 3193.814
$ws.Range("K31").Value = 12432.486
$ws.Range("L31").Value = 3193.814
$ws.Range("M31").Value = -12137.486
$ws.Range("N31").Value = -3783.814

$ws.Range("H34").Value = 7466.7
$ws.Range("I34").Value = 12432.486
$ws.Range("J34").Value = 3193.814
$ws.Range("K34").Value = 12432.486
$ws.Range("L34").Value = 3193.814
$ws.Range("M34").Value = -12230.486
$ws.Range("N34").Value = -3597.814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2500
$ws.Range("J55").Value = 2500
$ws.Range("L55").Value = 7500
$ws.Range("N55").Value = -7854

$ws.Range("H68").Value = 1846.9474
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 1866.2222
$ws.Range("K68").Value = 4500
$ws.Range("L68").Value = 5598.6666
$ws.Range("M68").Value = -3689
$ws.Range("N68").Value = -7220.6666

$ws.Range("H71").Value = 1846.9474
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 1866.2222
$ws.Range("K71").Value = 13500
$ws.Range("L71").Value = 16795.9998
$ws.Range("M71").Value = -9444
$ws.Range("N71").Value = -24907.9998

$ws.Range("H76").Value = 4566.7856
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 4764.231
$ws.Range("K76").Value = 6000
$ws.Range("L76").Value = 14292.693
$ws.Range("M76").Value = -5617
$ws.Range("N76").Value = -15058.693

$ws.Range("H79").Value = 4566.7856
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 4764.231
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 14292.693
$ws.Range("M79").Value = -4674
$ws.Range("N79").Value = -16944.693

$ws.Range("H107").Value = 21078
$ws.Range("I107").Value = 50195
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 150585
$ws.Range("L107").Value = 4999.9998
$ws.Range("M107").Value = -148665
$ws.Range("N107").Value = -8839.9998

$ws.Range("H131").Value = 125849.86
$ws.Range("J131").Value = 136000.94
$ws.Range("L131").Value = 408002.82
$ws.Range("N131").Value = -418082.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4802.96
$ws.Range("I126").Value = 3497.2856
$ws.Range("J126").Value = 6464.727
$ws.Range("K126").Value = 10491.8568
$ws.Range("L126").Value = 19394.181
$ws.Range("M126").Value = -8021.856800000001
$ws.Range("N126").Value = -24334.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4839.087
$ws.Range("I61").Value = 2517.6365
$ws.Range("J61").Value = 6967.0835
$ws.Range("K61").Value = 2517.6365
$ws.Range("L61").Value = 6967.0835
$ws.Range("M61").Value = -2315.6365
$ws.Range("N61").Value = -7371.0835

$ws.Range("H113").Value = 4839.087
$ws.Range("I113").Value = 2517.6365
$ws.Range("J113").Value = 6967.0835
$ws.Range("K113").Value = 2517.6365
$ws.Range("L113").Value = 6967.0835
$ws.Range("M113").Value = -347.6365000000001
$ws.Range("N113").Value = -11307.0835

$ws.Range("H132").Value = 1942.5385
$ws.Range("I132").Value = 1214.7142
$ws.Range("J132").Value = 4999.4
$ws.Range("K132").Value = 3644.1426
$ws.Range("L132").Value = 14998.2
$ws.Range("M132").Value = -1114.1426
$ws.Range("N132").Value = -20058.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 30271
$ws.Range("J68").Value = 30271
$ws.Range("L68").Value = 30271
$ws.Range("N68").Value = -31893

$ws.Range("H71").Value = 30271
$ws.Range("J71").Value = 30271
$ws.Range("L71").Value = 90813
$ws.Range("N71").Value = -98925

$ws.Range("H136").Value = 1167.4722
$ws.Range("I136").Value = 675.625
$ws.Range("J136").Value = 2151.1667
$ws.Range("K136").Value = 2026.875
$ws.Range("L136").Value = 6453.500100000001
$ws.Range("M136").Value = 523.125
$ws.Range("N136").Value = -11553.5001
